$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $range = $ws.Range($cellRef)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue 'D2' '26.428.62'
Set-TextValue 'E2' '  -3.76%  '
Set-TextValue 'D3' '1.770.18'
Set-TextValue 'E3' '  -3.01%  '
Set-TextValue 'D4' '1.003'
Set-TextValue 'E4' '  +0.09%  '
Set-TextValue 'D5' '1.002'
Set-TextValue 'E5' '  +0.05%  '
Set-TextValue 'D6' '306.33'
Set-TextValue 'E6' '  -2.07%  '
Set-TextValue 'D7' '0.4296'
Set-TextValue 'E7' '  +0.88%  '
Set-TextValue 'D8' '0.3660'
Set-TextValue 'E8' '  +1.49%  '
Set-TextValue 'D9' '0.07200'
Set-TextValue 'D10' '0.8501'
Set-TextValue 'E10' '  -1.20%  '
Set-TextValue 'D11' '20.36'
Set-TextValue 'E11' '  -0.98%  '
Set-TextValue 'D12' '1.765.89'
Set-TextValue 'E12' '  -4.42%  '
Set-TextValue 'D13' '6.429'
Set-TextValue 'E13' '  -0.63%  '
Set-TextValue 'D14' '5.235'
Set-TextValue 'E14' '  -2.68%  '
Set-TextValue 'D15' '0.06923'
Set-TextValue 'E15' '  -0.07%  '
Set-TextValue 'E16' '  -0.10%  '
Set-TextValue 'D17' '79.33'
Set-TextValue 'E17' '  -1.79%  '
Set-TextValue 'D18' '0.000008639'
Set-TextValue 'E18' '  -2.86%  '
Set-TextValue 'E19' '  +0.13%  '
Set-TextValue 'E20' '  -2.16%  '
Set-TextValue 'D21' '26.429.28'
Set-TextValue 'E21' '  -4.12%  '
Set-TextValue 'D22' '5.100'
Set-TextValue 'E22' '  -0.47%  '
Set-TextValue 'D23' '11.21'
Set-TextValue 'E23' '  +3.14%  '
Set-TextValue 'D24' '2.007.67'
Set-TextValue 'E24' '  -4.97%  '
Set-TextValue 'D25' '151.96'
Set-TextValue 'E25' '  -2.18%  '
Set-TextValue 'D26' '1.875'
Set-TextValue 'E26' '  -5.67%  '
Set-TextValue 'D27' '18.06'
Set-TextValue 'E27' '  -3.46%  '
Set-TextValue 'D28' '5.080'
Set-TextValue 'E28' '  -1.02%  '
Set-TextValue 'D29' '114.53'
Set-TextValue 'E29' '  +0.33%  '
Set-TextValue 'D30' '1.747'
Set-TextValue 'E30' '  -2.30%  '
Set-TextValue 'D31' '0.08951'
Set-TextValue 'E31' '  +0.61%  '
Set-TextValue 'D32' '0.7243'
Set-TextValue 'E32' '  -2.70%  '
Set-TextValue 'B33' 'Filecoin'
Set-TextValue 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue 'D33' '4.334'
Set-TextValue 'E33' '  -4.59%  '
Set-TextValue 'B34' 'ARBITRUM'
Set-TextValue 'C34' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D34' '1.112'
Set-TextValue 'E34' '  -0.44%  '
Set-TextValue 'E35' '  +0.05%  '
Set-TextValue 'D36' '2.744'
Set-TextValue 'E36' '  -8.15%  '
Set-TextValue 'D37' '1.081'
Set-TextValue 'E37' '  -0.30%  '
Set-TextValue 'D38' '0.05160'
Set-TextValue 'E38' '  -1.77%  '
Set-TextValue 'D39' '0.01888'
Set-TextValue 'E39' '  -1.73%  '
Set-TextValue 'D40' '0.4921'
Set-TextValue 'E40' '  -2.95%  '
Set-TextValue 'D41' '0.1608'
Set-TextValue 'E41' '  -2.59%  '
Set-TextValue 'D42' '2.573'
Set-TextValue 'E42' '  -7.69%  '
Set-TextValue 'D43' '6.268'
Set-TextValue 'E43' '  -1.51%  '
Set-TextValue 'D44' '8.007'
Set-TextValue 'E44' '  -3.87%  '
Set-TextValue 'D45' '104.87'
Set-TextValue 'E45' '  -1.44%  '
Set-TextValue 'B46' 'EnergySwap'
Set-TextValue 'C46' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '10.22'
Set-TextValue 'E46' '  -1.95%  '
Set-TextValue 'B47' 'PaxDollar'
Set-TextValue 'C47' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-TextValue 'D47' '1.002'
Set-TextValue 'E47' '  +0.07%  '
Set-TextValue 'D48' '0.4494'
Set-TextValue 'E48' '  -3.93%  '
Set-TextValue 'D49' '0.06193'
Set-TextValue 'E49' '  -4.04%  '
Set-TextValue 'D50' '1.589'
Set-TextValue 'E50' '  -1.41%  '
Set-TextValue 'D51' '1.741'
Set-TextValue 'E51' '  +2.77%  '
